$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values: switch from STM servo figures to ESP figures
$ws.Range("G13").Value = 251
$ws.Range("G15").Value = 440
$ws.Range("C21").Value = 251
$ws.Range("D21").Value = 440

# Update the active selection on the sheet view
$ws.Range("P15").Select()
